$wb = $excel.ActiveWorkbook

# zh-cn sheet: rows 21-22 both reference the same "Correspond Handoff
# Datetime" (D) and "Correspond Handback DateTime" (G) shared strings for
# the 2a54065f... handoff, so update every cell sharing that value.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D21").Value = "2016-03-04 11:25:20"
$wsZhCn.Range("D22").Value = "2016-03-04 11:25:20"
$wsZhCn.Range("G21").Value = "2016-03-04 11:26:15"
$wsZhCn.Range("G22").Value = "2016-03-04 11:26:15"

# de-de sheet: rows 10-11, same pattern.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D10").Value = "2016-03-04 11:25:33"
$wsDeDe.Range("D11").Value = "2016-03-04 11:25:33"
$wsDeDe.Range("G10").Value = "2016-03-04 11:26:41"
$wsDeDe.Range("G11").Value = "2016-03-04 11:26:41"
